# "andy text update august 27"
# Update the Skeena sheet's Region 6 Lake Babine Nation fisheries notice rows:
#  - Expand the abbreviated "Region 6-LBN" label to "Region 6-Lake Babine Nation"
#    on the existing Aug 15-21 notice (FN0821) row.
#  - Add a new fisheries notice row (FN0853, Aug 23-29) for the same region.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skeena")

# Append the new fisheries notice as row 37 (left-to-right so new shared
# strings are interned in the same order they first appear).
$ws.Range("A37").Value = "FN0853"
$ws.Range("B37").Value = "Aboriginal"
$ws.Range("C37").Value = "Aug 23-29"
$ws.Range("C37").NumberFormat = "d-mmm"
$ws.Range("D37").Value = "Sockeye"
$ws.Range("E37").Value = "Selective Gear"
$ws.Range("F37").Value = "Region 6-Lake Babine Nation"
$ws.Range("G37").Value = 7
$ws.Range("I37").Value = "Sockeye retention only"

# Expand the region name on the existing row (FN0821 / Aug 15-21), reusing
# the shared string just created above.
$ws.Range("F35").Value = "Region 6-Lake Babine Nation"

$ws.Range("I38").Select() | Out-Null
